# Generate Report for Handoff
# Appends a new localization-status row (for file
# 63d8a547-a403-4136-a001-3dc178b0b903.md) to the Overview sheet and to
# each per-locale sheet (zh-cn, de-de), mirroring the existing
# 60acfe56-... row already present on each sheet.
#
# NOTE: every text value is written with a leading apostrophe. That is
# Excel's "force text" quote-prefix, which keeps values such as "True" /
# "False" / date-look-alike strings stored as plain shared-string text
# (matching the source data) instead of being auto-coerced to Boolean /
# date cells. A lone apostrophe forces an (otherwise no-op) empty-string
# literal to still be written out as a real text cell.

$wb = $excel.ActiveWorkbook

$newId   = "63d8a547-a403-4136-a001-3dc178b0b903"
$newMd   = "$newId.md"
$commit  = "20cdb2bd7edc8d71514cb977f85b0e6fa6ee6fc8"
$repoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eefdf78a4dfe36567b77a2d78daf2e5f59b700a5/e2e/$newMd"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: "Overview" -> new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "'$newMd"
$wsOverview.Range("B3").Value = "'e2e\$newMd"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $repoUrl, "", "", "e2e\$newMd") | Out-Null
$wsOverview.Range("C3").Value = "'.md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "'Ready for handoff"
$wsOverview.Range("F3").Value = "'Ready for handoff"
$wsOverview.Range("G3").NumberFormat = $dateFmt
$wsOverview.Range("G3").Value = "'2016-09-07 04:52:52"

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn" -> new row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "'$newMd"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $repoUrl, "", "", $newMd) | Out-Null
$wsZhCn.Range("B3").Value = "'.md"
$wsZhCn.Range("C3").Value = "'Ready for handoff"
$wsZhCn.Range("D3").Value = "'e2e"
$wsZhCn.Range("E3").Value = "'ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "'$newId.$commit.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("H3").Value = "'2016-09-07 04:52:46"
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

# ---------------------------------------------------------------------
# Sheet 3: "de-de" -> new row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "'$newMd"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $repoUrl, "", "", $newMd) | Out-Null
$wsDeDe.Range("B3").Value = "'.md"
$wsDeDe.Range("C3").Value = "'Ready for handoff"
$wsDeDe.Range("D3").Value = "'e2e"
$wsDeDe.Range("E3").Value = "'ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "'$newId.$commit.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("H3").Value = "'2016-09-07 04:52:52"
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

Write-Output "Added handoff row for $newMd to Overview, zh-cn and de-de sheets."
